$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping CSV2openEHR")

$ws.Range("B3").Value = "bericht/context/status"
$ws.Range("B4").Value = "bericht/context/umg_personenname<<index>>/name_strukturiert/art_des_pr_fix"
$ws.Range("C4").Value = 0
$ws.Range("B5").Value = "bericht/context/umg_personenname<<index>>/name_strukturiert/vorname"
$ws.Range("C5").Value = 0

$ws.Range("C12").Select()
